# "alterações feitas de modo a ficarem mais corretas"
# Fix mixed-up / duplicated "Meat.*" ViewModel names on the Folha1 sheet:
#   - H4  (next to "Info")     was "Meat.info"        -> should be "Meat.moreInfo"
#   - H7  (next to "See more") was "Meat.moreInfo"     -> should be blank
#   - D37 (next to "Meat details") was "Meat.generalInfo" -> should be "Meat.details"
#   - D38 (next to "Expand")   was "Meat.moreInfo"     -> should be blank

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H4").Value = "Meat.moreInfo"
$ws.Range("H7").Value = ""

$ws.Range("D37").Value = "Meat.details"
$ws.Range("D38").Value = ""

# Update the view: scroll so column D is the leftmost visible column,
# zoom to 130%, and move the selection to G13.
$win = $excel.ActiveWindow
$win.Zoom = 130
$win.ScrollColumn = 4
$win.ScrollRow = 1

$ws.Range("G13").Select()
